$wb = $excel.ActiveWorkbook

# Overview sheet: update "Latest HO Xliff Generate Date" for first row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-28 01:02:17"

# zh-cn sheet: update Correspond Handoff Datetime and Correspond Handback DateTime for first row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-28 01:02:12"
$wsZhCn.Range("K2").Value = "2016-08-28 01:02:28"

# de-de sheet: update Correspond Handback DateTime for first row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-28 01:02:34"
